$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new numeric-looking text must be forced to remain as text
# (matching the source workbook, which stores these as inline strings).
$textCells = @("D5", "D6", "D8", "D11", "D17", "D18", "D19", "D25", "D27", "D30", "D37", "D42", "D43", "D46", "D48", "D51")
foreach ($c in $textCells) {
  $ws.Range($c).NumberFormat = "@"
}

# Apply updated Price (D) and Volume(1h) (E) values row by row.
$ws.Range("D2").Value = "26.899.25"
$ws.Range("E2").Value = "  -1.35%  "
$ws.Range("D3").Value = "1.564.44"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "206.03"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").Value = "0.489"
$ws.Range("E6").Value = "  -1.31%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").Value = "21.77"
$ws.Range("E8").Value = "  -1.78%  "
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("E10").Value = "  -1.31%  "
$ws.Range("D11").Value = "0.0866"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "1.563.09"
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("E14").Value = "  -1.45%  "
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("D16").Value = "26.867.22"
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("D17").Value = "61.22"
$ws.Range("E17").Value = "  -3.16%  "
$ws.Range("D18").Value = "7.37"
$ws.Range("E18").Value = "  +1.78%  "
$ws.Range("D19").Value = "214.55"
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("E20").Value = "  -1.46%  "
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("E23").Value = "  -2.71%  "
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").Value = "153.99"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("D27").Value = "14.95"
$ws.Range("E27").Value = "  +0.32%  "
$ws.Range("E28").Value = "  -0.29%  "
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("D30").Value = "0.0464"
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("E31").Value = "  -3.31%  "
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").Value = "1.404.04"
$ws.Range("E33").Value = "  +1.43%  "
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("E35").Value = "  -1.50%  "
$ws.Range("E36").Value = "  -1.18%  "
$ws.Range("D37").Value = "0.920"
$ws.Range("E37").Value = "  -2.46%  "
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("D42").Value = "0.997"
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("D43").Value = "5.34"
$ws.Range("E43").Value = "  +2.28%  "
$ws.Range("E44").Value = "  -2.21%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "63.17"
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("D47").Value = "1.699.03"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "86.22"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("E49").Value = "  +2.40%  "
$ws.Range("D50").Value = "0.0₇0981"
$ws.Range("E50").Value = "  -1.41%  "
$ws.Range("D51").Value = "0.0946"
$ws.Range("E51").Value = "  +0.17%  "

# Restore the Normal style on the forced-text cells so no stray style index
# is introduced (the source file keeps these cells styleless).
foreach ($c in $textCells) {
  $ws.Range($c).Style = "Normal"
}

Write-Host "Updated cryptos list values"
